$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

for ($row = 2; $row -le 22; $row++) {
    $ws.Cells.Item($row, 5).Value = "fullRNASeq"  # Column E
}
